$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-5 get cyclically rotated:
#   new row2 = old row4
#   new row3 = old row2
#   new row4 = old row5
#   new row5 = old row3
# Capture the original values for the columns that change (D, M, N, O, P, R, S)
# before overwriting any of them.

$cols = @("D","M","N","O","P","R","S")

$orig = @{}
foreach ($r in 2..5) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $orig[$r] = $rowVals
}

$mapping = @{ 2 = 4; 3 = 2; 4 = 5; 5 = 3 }

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    foreach ($col in $cols) {
        $ws.Range("$col$newRow").Value2 = $orig[$oldRow][$col]
    }
}
